$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bus voltage magnitude results for the 380 kV case (Case_5_240)
# Column B changes uniformly from 1.05 -> 1.02 (slack/reference voltage setpoint),
# columns C-F and I-M are recomputed power-flow results for rows 2-25 (bus indices 0-23).

$colB = @(1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colB[$i]
}

$colC = @(1.061376949907918, 1.063115506579163, 1.064237477795957, 1.064708453617011, 1.064787491869683, 1.064243773731071, 1.061965131824149, 1.057926278518517, 1.05521688191818, 1.054039476293978, 1.05360148287724, 1.053695463783195, 1.054003285042077, 1.054192856847983, 1.05529493197204, 1.055985091907281, 1.056387244719838, 1.056524300167528, 1.055911086386575, 1.053912657518111, 1.05265238114108, 1.053320842545805, 1.055944527549458, 1.058973313626455)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $colC[$i]
}

$colD = @(1.059268964046522, 1.060615663001646, 1.061484340803692, 1.061848891014127, 1.06191006326648, 1.061489214448988, 1.059724658319982, 1.056593916734176, 1.054491703971969, 1.05357769891199, 1.053237622533506, 1.053310596271111, 1.053549599883853, 1.053696781408041, 1.054552283690237, 1.055087908847166, 1.055399970663223, 1.055506315140134, 1.0550304785949, 1.053479235235767, 1.052500578072912, 1.053019702792322, 1.055056429963859, 1.057405890747683)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $colD[$i]
}

$colE = @(1.066137741294542, 1.067704699373687, 1.068715636053693, 1.069139931483659, 1.069211131755969, 1.068721308251306, 1.06666792934176, 1.063026107187637, 1.060581569445362, 1.059518916211598, 1.059123559119494, 1.059208393803262, 1.059486249030653, 1.059657359453906, 1.060652005350228, 1.061274797149992, 1.061637661563657, 1.061761321685708, 1.061208018926082, 1.059404445402539, 1.058266753853074, 1.058870223216249, 1.061238194365253, 1.063970478271023)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $colE[$i]
}

$colF = @(1.075198332236494, 1.076849878206589, 1.077915825873264, 1.078363314474427, 1.078438412891997, 1.077921807717315, 1.075757049959421, 1.071921074080035, 1.069348503321495, 1.068230739508575, 1.067814960944811, 1.067904173888959, 1.068196383283604, 1.068376344273402, 1.069422603629957, 1.070077858936848, 1.070459690487784, 1.070589823303555, 1.070007594406864, 1.068110351359475, 1.066914049699699, 1.067548562611164, 1.070039345059675, 1.072915390122201)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $colF[$i]
}

$colI = @(1.046764804841948, 1.047272688990757, 1.047598461532078, 1.047734736153023, 1.047757577524757, 1.047600285105713, 1.04693704246, 1.04574616961598, 1.044937032602859, 1.044582983561595, 1.044450913589352, 1.044479268509471, 1.044572078085721, 1.044629186677987, 1.044960451437829, 1.045167253356162, 1.045287522075185, 1.045328470520918, 1.045145102268322, 1.044544763500219, 1.044164061188242, 1.044366188483059, 1.045155112489346, 1.046056697567916)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item(2 + $i, 9).Value = $colI[$i]
}

$colJ = @(1.066352946779296, 1.067741759823263, 1.068637107034582, 1.06901273142618, 1.069075755054717, 1.068642129192429, 1.066822994961913, 1.063591560669158, 1.061419086660913, 1.060473893723973, 1.060122116153484, 1.060197605143117, 1.060444829861253, 1.060597061171643, 1.061481719959844, 1.062035428334256, 1.062357964164015, 1.062467867575769, 1.061976065551843, 1.060372047536507, 1.059359533188601, 1.059896670892451, 1.062002890375388, 1.064430111371474)
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $ws.Cells.Item(2 + $i, 10).Value = $colJ[$i]
}

$colK = @(1.061998929533475, 1.063158477468311, 1.063905521738347, 1.064218808669961, 1.064271366036323, 1.06390971090716, 1.062391485495858, 1.05969074625806, 1.057872509297541, 1.057080839238979, 1.056786109889671, 1.056849360760371, 1.057056490516094, 1.057184021120405, 1.057924956870776, 1.058388549771934, 1.058658535708231, 1.058750523081292, 1.05833885417022, 1.056995514503065, 1.056147033275773, 1.056597200323916, 1.058361310771012, 1.060392031867537)
for ($i = 0; $i -lt $colK.Length; $i++) {
    $ws.Cells.Item(2 + $i, 11).Value = $colK[$i]
}

$colL = @(1.068849097986149, 1.070229715825728, 1.071119587538366, 1.071492868959049, 1.071555496791676, 1.071124578550352, 1.069316412467505, 1.066102971527101, 1.063941580305755, 1.063000966382796, 1.062650856330957, 1.06272598914517, 1.06297204103932, 1.063123545308723, 1.064003905052445, 1.064554857667895, 1.064875764911992, 1.06498510943145, 1.064495792747075, 1.062899605004452, 1.061891820872512, 1.062426469389806, 1.064522483045269, 1.066937030238204)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $ws.Cells.Item(2 + $i, 12).Value = $colL[$i]
}

$colM = @(1.07788553521846, 1.079352306830526, 1.080298215184435, 1.08069512389279, 1.080761722931379, 1.080303521632084, 1.078381905973982, 1.074970759497738, 1.072679035291621, 1.071682331557867, 1.071311438857026, 1.071391027289082, 1.071651687262974, 1.071812198882519, 1.072745089785632, 1.073329086053814, 1.07366930042944, 1.073785233842666, 1.073266472416411, 1.071574948090332, 1.070507520495799, 1.071073758942466, 1.073294766123911, 1.075855669620246)
for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item(2 + $i, 13).Value = $colM[$i]
}
